$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 14 with ID, Description and Hierarchy values
$ws.Range("A14").Value = "ASSUM"
$ws.Range("B14").Value = "Salary Assumptions"
$ws.Range("C14").Value = "<root>"

# Update the active selection to match the new last-edited cell
$ws.Range("C14").Select()
